$d = $word.ActiveDocument

# --- Change 1: merge split runs "THU Dec 07" + " 09:57:40 PST 2017" into one run ---
$d.Content.Find.Execute("THU Dec 07 09:57:40 PST 2017", $false, $false, $false, $false, $false, $true, 1, $false, "THU Dec 07 09:57:40 PST 2017", 2) | Out-Null

# --- Change 2: append a new purchase-details block after the last "CASH AND CLEARD" paragraph ---

# locate the paragraph that holds the final "- CASH AND CLEARD" entry
$paras = $d.Paragraphs
$n = $paras.Count
$lastIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*CASH AND CLEARD*") {
        $lastIdx = $i
    }
}

$insPara = $d.Paragraphs.Item($lastIdx + 1)
$insPoint = $insPara.Range
$insPoint.Collapse(1)

$block = ""
$block += "`r"
$block += "THU Dec 14`r 10:16:53 PST 2017`r"
$block += "Person Name`t`t`t`t- TSV`r"
$block += "Bill number`t`t`t`t- 1936`r"
$block += "---------------------------------------------------------------`r"
$block += "Item Name`t`t`t`t- CHOWCHOW`r"
$block += "Number of Pockets`t`t`t- 2`r"
$block += "Number of KGs`t`t`t- 139`r"
$block += "Rate`t`t`t`t`t- 6`r"
$block += "Total Price`t`t`t`t- 834.0`r"
$block += "Amount balance`t`t`t- 834.0`r"
$block += "`r"
$block += "`r"

$insPoint.InsertBefore($block)

# fix up the date paragraph: split "THU Dec 14" / " 10:16:53 PST 2017" into two
# separate runs (identical formatting) by joining two paragraphs instead of typing
# one after another (which the app would otherwise coalesce into a single run).
$dateP = $d.Paragraphs.Item($lastIdx + 2)
$markStart = $dateP.Range.End - 1
$d.Range($markStart, $markStart + 1).Delete()

# bold the "Amount balance" line (paragraph + paragraph mark)
$boldPara = $d.Paragraphs.Item($lastIdx + 10)
$boldPara.Range.Bold = 1

Write-Host "done"
